# İş Takip Güncellemesi - 26.11.2025 10:42:54
# Shifts every date-like text value (yyyy-MM-dd) in the tracked columns back by one day.

$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-MM-dd"
$culture = [System.Globalization.CultureInfo]::InvariantCulture

function Shift-CellDateBack {
    param($ws, [int]$row, [int]$col, [bool]$skip)
    if ($skip) {
        return
    }
    $cell = $ws.Cells.Item($row, $col)
    $val = $cell.Value()
    if ($val -match '^\d{4}-\d{2}-\d{2}$') {
        $dt = [DateTime]::ParseExact($val, $dateFormat, $culture)
        $newVal = $dt.AddDays(-1).ToString($dateFormat, $culture)
        # Force text storage so Excel doesn't coerce the string into a date serial value.
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
    }
}

# Sheet 1: "İş Takip Listesi" -> columns J (İŞE BAŞLAMA/YER TESLİMİ) and K (İHALE BİTİŞ TARİHİ)
$ws1 = $wb.Worksheets.Item(1)
$usedRows1 = $ws1.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows1; $r++) {
    Shift-CellDateBack $ws1 $r 10 $false   # J
    Shift-CellDateBack $ws1 $r 11 $false   # K
}

# Sheet 2: "Güncelleme" -> columns I, J, N, P
# Note: row 14's N column (BİLGİLENDİRME İLANI TARİHİ) is intentionally left as-is.
$ws2 = $wb.Worksheets.Item(2)
$usedRows2 = $ws2.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows2; $r++) {
    $skipN = $false
    if ($r -eq 14) {
        $skipN = $true
    }
    Shift-CellDateBack $ws2 $r 9  $false   # I
    Shift-CellDateBack $ws2 $r 10 $false   # J
    Shift-CellDateBack $ws2 $r 14 $skipN   # N (skip row 14)
    Shift-CellDateBack $ws2 $r 16 $false   # P
}
